# Trade #45 closed at 2026-02-18 00:19:14 - unknown UNKNOWN +0.000%
#
# This script:
#  1. Updates the Summary sheet roll-up metrics.
#  2. Updates the Strategy Status row for EMAArbitrage.
#  3. Closes trade #74 (EMAArbitrage) on the "All Trades" sheet and its
#     per-strategy "EMAArbitrage" sheet (status OPEN -> CLOSED, fills in
#     exit price / P&L / capital / exit reason / duration).
#  4. Appends two brand-new OPEN trades (#103 momentum, #104
#     HighProbConvergence) to "All Trades" and to their respective
#     per-strategy sheets.

$wb = $excel.ActiveWorkbook

# Helper: force a date/time-looking literal to be written as plain text
# instead of being auto-converted to a date serial number by Excel's
# input parser.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.44   # Current Capital
$summary.Range("B4").Value = 0.54      # Total P&L $
$summary.Range("B5").Value = 0.15      # Total P&L %
$summary.Range("B6").Value = 73        # Total Trades
$summary.Range("B7").Value = 37        # Winning Trades
$summary.Range("B9").Value = 50.68     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - EMAArbitrage row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C2").Value = 100.29     # Capital
$status.Range("D2").Value = 4          # Trades
$status.Range("E2").Value = 0.29       # P&L $
$status.Range("F2").Value = 0.29       # P&L %
$status.Range("G2").Value = 75         # Win Rate %

# ---------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close trade #74 (EMAArbitrage) on row 75
$allTrades.Range("G75").Value = 0.99        # Exit Price
$allTrades.Range("H75").Value = "CLOSED"    # Status
$allTrades.Range("I75").Value = 2.0619      # P&L %
$allTrades.Range("J75").Value = 0.02        # P&L $
$allTrades.Range("K75").Value = 100.29      # Capital After
$allTrades.Range("L75").Value = "early_exit" # Exit Reason
$allTrades.Range("M75").Value = 0.12        # Duration (min)

# New trade #103 (momentum) -> row 104
$allTrades.Range("A104").Value = 103
Set-TextValue $allTrades.Range("B104") "2026-02-18"
Set-TextValue $allTrades.Range("C104") "00:19:07"
$allTrades.Range("D104").Value = "momentum"
$allTrades.Range("E104").Value = "UP"
$allTrades.Range("F104").Value = 0.97
$allTrades.Range("H104").Value = "OPEN"
$allTrades.Range("I104").Value = 0
$allTrades.Range("J104").Value = 0
$allTrades.Range("K104").Value = 99.6787371310913
$allTrades.Range("M104").Value = 0
$allTrades.Range("N104").Value = 0
$allTrades.Range("O104").Value = 0
$allTrades.Range("P104").Value = 0.9
$allTrades.Range("Q104").Value = "Upward momentum: 3.960% over 10 samples"

# New trade #104 (HighProbConvergence) -> row 105
$allTrades.Range("A105").Value = 104
Set-TextValue $allTrades.Range("B105") "2026-02-18"
Set-TextValue $allTrades.Range("C105") "00:19:09"
$allTrades.Range("D105").Value = "HighProbConvergence"
$allTrades.Range("E105").Value = "DOWN"
$allTrades.Range("F105").Value = 0.03
$allTrades.Range("H105").Value = "OPEN"
$allTrades.Range("I105").Value = 0
$allTrades.Range("J105").Value = 0
$allTrades.Range("K105").Value = 100.0565626577805
$allTrades.Range("M105").Value = 0
$allTrades.Range("N105").Value = 0
$allTrades.Range("O105").Value = 0
$allTrades.Range("P105").Value = 0.95
$allTrades.Range("Q105").Value = "Mean reversion DOWN: price 3.75% above mean (z=4.36)"

# ---------------------------------------------------------------------
# 4. momentum sheet - append trade #103 as row 24
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A24").Value = 103
Set-TextValue $momentum.Range("B24") "2026-02-18"
Set-TextValue $momentum.Range("C24") "00:19:07"
$momentum.Range("D24").Value = "momentum"
$momentum.Range("E24").Value = "UP"
$momentum.Range("F24").Value = 0.97
$momentum.Range("H24").Value = "OPEN"
$momentum.Range("I24").Value = 0
$momentum.Range("J24").Value = 0
$momentum.Range("K24").Value = 99.6787371310913
$momentum.Range("L24").Value = 0
$momentum.Range("M24").Value = 0
$momentum.Range("N24").Value = 0.9
$momentum.Range("O24").Value = "Upward momentum: 3.960% over 10 samples"
$momentum.Range("Q24").Value = 0

# ---------------------------------------------------------------------
# 5. HighProbConvergence sheet - append trade #104 as row 12
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("A12").Value = 104
Set-TextValue $hpc.Range("B12") "2026-02-18"
Set-TextValue $hpc.Range("C12") "00:19:09"
$hpc.Range("D12").Value = "HighProbConvergence"
$hpc.Range("E12").Value = "DOWN"
$hpc.Range("F12").Value = 0.03
$hpc.Range("H12").Value = "OPEN"
$hpc.Range("I12").Value = 0
$hpc.Range("J12").Value = 0
$hpc.Range("K12").Value = 100.0565626577805
$hpc.Range("L12").Value = 0
$hpc.Range("M12").Value = 0
$hpc.Range("N12").Value = 0.95
$hpc.Range("O12").Value = "Mean reversion DOWN: price 3.75% above mean (z=4.36)"
$hpc.Range("Q12").Value = 0

# ---------------------------------------------------------------------
# 6. EMAArbitrage sheet - close trade #74 on row 5
# ---------------------------------------------------------------------
$ema = $wb.Worksheets.Item("EMAArbitrage")
$ema.Range("G5").Value = 0.99        # Exit Price
$ema.Range("H5").Value = "CLOSED"    # Status
$ema.Range("I5").Value = 2.0619      # P&L %
$ema.Range("J5").Value = 0.02        # P&L $
$ema.Range("K5").Value = 100.29      # Capital After
$ema.Range("P5").Value = "early_exit" # Exit Reason
$ema.Range("Q5").Value = 0.12        # Duration (min)
